$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 52: GenomeWeb link for the new Burning Rock story ---
$ws.Range("A52").Value = "https://www.genomeweb.com/cancer/burning-rocks-sequencing-test-approved-japan-cdx-truqap-breast-cancer"
$ws.Range("B52").Value = "CDx"
$ws.Range("C52").Value = "Burning Rock's Sequencing Test Approved in Japan as a CDx for Truqap in Breast Cancer"

# --- Row 53: 360Dx link for the same story ---
$ws.Range("A53").Value = "https://www.360dx.com/cancer/burning-rocks-sequencing-test-approved-japan-cdx-truqap-breast-cancer"
$ws.Range("B53").Value = "CDx"
$ws.Range("C53").Value = "Burning Rock's Sequencing Test Approved in Japan as a CDx for Truqap in Breast Cancer"

# Wire up the hyperlinks (column A) for the new rows, same as the existing link rows
$ws.Hyperlinks.Add($ws.Range("A52"), "https://www.genomeweb.com/cancer/burning-rocks-sequencing-test-approved-japan-cdx-truqap-breast-cancer")
$ws.Hyperlinks.Add($ws.Range("A53"), "https://www.360dx.com/cancer/burning-rocks-sequencing-test-approved-japan-cdx-truqap-breast-cancer")

# Re-apply the same formatting used by the other "link" cells in column A
# (Hyperlinks.Add creates its own style entry; copy the existing look instead)
$ws.Range("A2").Copy()
$ws.Range("A52").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("A53").PasteSpecial(-4122)
